$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 19 (pushes existing rows 19..121 down to 20..122)
$ws.Rows.Item(19).Insert()

# Populate the newly inserted row 19 with its new record
$ws.Cells.Item(19, 1).Value = 8
$ws.Cells.Item(19, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(19, 3).Value = "Coquimbo"
$ws.Cells.Item(19, 4).Value = 44831
$ws.Cells.Item(19, 5).Value = 4
$ws.Cells.Item(19, 6).Value = 100112052
$ws.Cells.Item(19, 7).Value = "Albahaca"
$ws.Cells.Item(19, 8).Value = "Sin especificar"
$ws.Cells.Item(19, 9).Value = "Primera"
$ws.Cells.Item(19, 10).Value = 1200
$ws.Cells.Item(19, 11).Value = 4000
$ws.Cells.Item(19, 12).Value = 4500
$ws.Cells.Item(19, 13).Value = 4250
$ws.Cells.Item(19, 14).Value = "`$/paquete"
$ws.Cells.Item(19, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(19, 16).Value = 4250
$ws.Cells.Item(19, 17).Value = 1
$ws.Cells.Item(19, 18).Value = "Hortaliza"

# Append a brand new record as row 123 (after the old row 121 shifted to 122)
$ws.Cells.Item(123, 1).Value = 8
$ws.Cells.Item(123, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(123, 3).Value = "Coquimbo"
$ws.Cells.Item(123, 4).Value = 44832
$ws.Cells.Item(123, 5).Value = 4
$ws.Cells.Item(123, 6).Value = 100112052
$ws.Cells.Item(123, 7).Value = "Albahaca"
$ws.Cells.Item(123, 8).Value = "Sin especificar"
$ws.Cells.Item(123, 9).Value = "Primera"
$ws.Cells.Item(123, 10).Value = 1400
$ws.Cells.Item(123, 11).Value = 4000
$ws.Cells.Item(123, 12).Value = 4500
$ws.Cells.Item(123, 13).Value = 4250
$ws.Cells.Item(123, 14).Value = "`$/paquete"
$ws.Cells.Item(123, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(123, 16).Value = 4250
$ws.Cells.Item(123, 17).Value = 1
$ws.Cells.Item(123, 18).Value = "Hortaliza"

# Apply the same date number-format used by every other "Fecha" cell in column D
$ws.Cells.Item(123, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
